$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $replace, 2) | Out-Null
}

Replace-Text "2025-10-19 Sunday" "2025-10-20 Monday"
Replace-Text "63×84=" "24×24="
Replace-Text "92×91=" "44×87="
Replace-Text "19×94=" "85×76="
Replace-Text "21×64=" "61×45="
Replace-Text "56×57=" "26×58="
Replace-Text "99×15=" "26×51="
Replace-Text "63×82=" "45×31="
Replace-Text "37×91=" "40×99="
Replace-Text "12×20=" "35×17="
Replace-Text "49×58=" "11×32="
Replace-Text "75×43=" "36×83="
Replace-Text "81×51=" "24×49="
Replace-Text "77×59=" "58×92="
Replace-Text "64×14=" "49×60="
Replace-Text "37×94=" "85×57="
Replace-Text "86×60=" "15×19="
Replace-Text "26×80=" "35×22="
Replace-Text "62×28=" "99×95="
Replace-Text "91×96=" "54×97="
Replace-Text "28×64=" "51×94="
Replace-Text "34×55=" "41×46="
Replace-Text "61×19=" "50×92="
Replace-Text "68×13=" "80×74="
Replace-Text "80×85=" "29×39="
Replace-Text "24×45=" "91×97="
